# Codebloggs Home/PWA HowToImprove.docx -- "lighthouse done + pull-request version"
#
# Changes applied (per the commit diff):
#  1. Add <w:lang w:val="en-CA"/> to the "Current note : ..." paragraph's
#     paragraph-mark rPr, to the "Current note : " run, and to the
#     "voir par titre" run.
#  2. Translate "voir par titre" -> "see per title", typed/split as two
#     runs ("see per t" + "itle") with identical (now en-CA) formatting.
#  3. Give the blank paragraph right after the table a pPr/rPr with
#     <w:lang w:val="en-CA"/> (it was a bare <w:p/>).
#  4. Merge the three separate "4", "/", "5" runs of "PWA OPTIMIZED : 4/5"
#     into a single "4/5" run.
#
# We do this by pulling the whole body as WordOpenXML, editing the raw
# markup with plain string surgery (so we get exact control over run
# boundaries / rPr contents), and feeding it back with InsertXML.

$d = $word.ActiveDocument

$xml = $d.Content.WordOpenXML

# --- 1a. paragraph-mark rPr for the "Current note :" paragraph ---------
$old = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="0F4761" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>'
$new = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="0F4761" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-CA"/></w:rPr></w:pPr>'
$xml = $xml.Replace($old, $new)

# --- 1b. "Current note : " run rPr --------------------------------------
$old = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="0F4761" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Current note : </w:t>'
$new = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="0F4761" w:themeColor="accent1" w:themeShade="BF"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve">Current note : </w:t>'
$xml = $xml.Replace($old, $new)

# --- 1c + 2. "voir par titre" run -> two en-CA runs "see per t" / "itle" -
$old = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="00B0F0"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>voir par titre</w:t>'
$new = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="00B0F0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-CA"/></w:rPr><w:t>see per t</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="00B0F0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-CA"/></w:rPr><w:t>itle</w:t>'
$xml = $xml.Replace($old, $new)

# --- 3. blank paragraph right after the table ---------------------------
$old = '</w:tbl><w:p w14:paraId="79FC654F" w14:textId="77777777" w:rsidR="007515E2" w:rsidRPr="00063155" w:rsidRDefault="007515E2" w:rsidP="007515E2"/>'
$new = '</w:tbl><w:p w14:paraId="79FC654F" w14:textId="77777777" w:rsidR="007515E2" w:rsidRPr="00063155" w:rsidRDefault="007515E2" w:rsidP="007515E2"><w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr></w:p>'
$xml = $xml.Replace($old, $new)

# --- 4. merge "4" + "/" + "5" runs into a single "4/5" run --------------
$old = '<w:r w:rsidRPr="00063155"><w:rPr><w:b/><w:bCs/><w:color w:val="FFC000"/><w:lang w:val="en-CA"/></w:rPr><w:t>4</w:t></w:r><w:r w:rsidRPr="00063155"><w:rPr><w:b/><w:bCs/><w:color w:val="FFC000"/><w:lang w:val="en-CA"/></w:rPr><w:t>/</w:t></w:r><w:r w:rsidRPr="00063155"><w:rPr><w:b/><w:bCs/><w:color w:val="FFC000"/><w:lang w:val="en-CA"/></w:rPr><w:t>5</w:t></w:r>'
$new = '<w:r w:rsidRPr="00063155"><w:rPr><w:b/><w:bCs/><w:color w:val="FFC000"/><w:lang w:val="en-CA"/></w:rPr><w:t>4/5</w:t></w:r>'
$xml = $xml.Replace($old, $new)

$d.Content.InsertXML($xml)
